# Update "想去人数" (want-to-go count) figures in column F across the
# workbook's sheets, matching the refreshed data snapshot
# ("Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 680
$ws.Range("F6").Value = 110
$ws.Range("F7").Value = 1150
$ws.Range("F9").Value = 49
$ws.Range("F10").Value = 2032
$ws.Range("F12").Value = 38
$ws.Range("F13").Value = 31
$ws.Range("F16").Value = 1478
$ws.Range("F18").Value = 551
$ws.Range("F19").Value = 388
$ws.Range("F20").Value = 388
$ws.Range("F21").Value = 720
$ws.Range("F22").Value = 449
$ws.Range("F23").Value = 2816
$ws.Range("F26").Value = 3205
$ws.Range("F27").Value = 658
$ws.Range("F28").Value = 526
$ws.Range("F29").Value = 231
$ws.Range("F30").Value = 976
$ws.Range("F31").Value = 730
$ws.Range("F33").Value = 690
$ws.Range("F34").Value = 668

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 84
$ws.Range("F10").Value = 63
$ws.Range("F13").Value = 78
$ws.Range("F20").Value = 96
$ws.Range("F21").Value = 183
$ws.Range("F22").Value = 131

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 376
$ws.Range("F6").Value = 386

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 376
$ws.Range("F8").Value = 680
$ws.Range("F9").Value = 110
$ws.Range("F11").Value = 1150
$ws.Range("F12").Value = 84
$ws.Range("F14").Value = 49
$ws.Range("F15").Value = 386
$ws.Range("F16").Value = 2032
$ws.Range("F17").Value = 2032
$ws.Range("F18").Value = 38
$ws.Range("F19").Value = 31
$ws.Range("F22").Value = 63
$ws.Range("F26").Value = 78
$ws.Range("F27").Value = 1478
$ws.Range("F28").Value = 1478
$ws.Range("F31").Value = 551
$ws.Range("F32").Value = 388
$ws.Range("F33").Value = 388
$ws.Range("F35").Value = 720
$ws.Range("F36").Value = 449
$ws.Range("F38").Value = 2816
$ws.Range("F40").Value = 3205
$ws.Range("F41").Value = 658
$ws.Range("F42").Value = 526
$ws.Range("F43").Value = 231
$ws.Range("F44").Value = 976
$ws.Range("F46").Value = 96
$ws.Range("F47").Value = 131
$ws.Range("F49").Value = 730
$ws.Range("F51").Value = 691
$ws.Range("F52").Value = 668
